# Fix Training Data Issue (#48)
# The "Date" column (BF) stored the literal source-file name
# ("5-26-2011-12") instead of the real game date. Because NBA box-score
# stats for a night get published under the next calendar day, the data
# was effectively off by one day. Correct value for this sheet: 2012-05-26.
# Every data row (2-31) gets its BF cell updated from "5-26-2011-12" to
# "2012-05-26".
#
# Implementation note: assigning a date-shaped literal like "2012-05-26"
# straight to Range.Value makes Excel auto-detect it as a date and convert
# it to a date serial number (pulling in a new number-format style along
# the way) instead of keeping it as literal text. To preserve the original
# plain-text cell content/type and leave cell styling untouched, the new
# text is first produced as a formula result ( ="2012-05-26" ) on a
# throwaway scratch cell, then copied over onto the real cell with
# PasteSpecial values-only. Pasting an already-resolved text value this
# way does not re-trigger Excel's "looks like a date" auto-conversion.
# The scratch row is deleted again afterwards so it leaves no trace.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 31
$oldDate = "5-26-2011-12"
$newDate = "2012-05-26"

# A scratch row far below the sheet's real data (1-31) so the temporary
# work can't collide with anything; it is removed again at the end.
$scratchRow = 100
$scratch = $ws.Range("ZZ$scratchRow")
$scratch.Formula = '="' + $newDate + '"'
$scratch.Copy()

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $target = $ws.Range("BF$r")
    if ($target.Value() -eq $oldDate) {
        $target.PasteSpecial(-4163)  # xlPasteValues - literal text, no format/formula carried over
    }
}

# Clean up the scratch row completely (contents + the row itself) so the
# sheet's used range/dimension is unaffected.
$ws.Range("ZZ$scratchRow").EntireRow.Delete()
